$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "面積（平方公尺）"
$ws.Range("D1").Value = "權利範圍(持分）"
$ws.Range("E1").Value = "所有權人"
$ws.Range("F1").Value = "登記（取得）時間"
$ws.Range("G1").Value = "登記（取得）原因"
$ws.Range("H1").Value = "取得價額"
$ws.Range("I1").Value = "property_category"
$ws.Range("J1").Value = "category"
$ws.Range("K1").Value = "date"
$ws.Range("L1").Value = "legislator_name"
$ws.Range("M1").Value = "legislator_id"
$ws.Range("N1").Value = "source_file"
$ws.Range("O1").Value = "index"
$ws.Range("A2").Value = 13
$ws.Range("B2").Value = "高雄市鼓山區龍北段00070000地號"
$ws.Range("C2").Value = 2275.44
$ws.Range("D2").Value = "10000分之110"
$ws.Range("E2").Value = "吳美惠"
$ws.Range("F2").Value = "102年07月22日"
$ws.Range("G2").Value = "買賣"
$ws.Range("H2").Value = "15150000(含2車位）"
$ws.Range("I2").Value = "land"
$ws.Range("J2").Value = "normal"
$ws.Range("K2").Value = "2013-12-24"
$ws.Range("L2").Value = "林國正"
$ws.Range("M2").Value = 1742
$ws.Range("N2").Value = "tmp399c1"
$ws.Range("O2").Value = 13

$ws = $wb.Worksheets.Item(2)
$ws.Range("B1").Value = "建物標示"
$ws.Range("C1").Value = "面積（平方公尺）"
$ws.Range("D1").Value = "權利範圍(持分）"
$ws.Range("E1").Value = "所有權人"
$ws.Range("F1").Value = "登記（取得）時間"
$ws.Range("G1").Value = "登記（取得）原因"
$ws.Range("H1").Value = "取得價額"
$ws.Range("A2").Value = 18
$ws.Range("B2").Value = "高雄市鼓山區龍北段00260000建號"
$ws.Range("C2").Value = 150.11
$ws.Range("D2").Value = "全部"
$ws.Range("E2").Value = "吳美惠"
$ws.Range("F2").Value = "102年07月22日"
$ws.Range("G2").Value = "買賣"
$ws.Range("H2").Value = 5480000
$ws.Range("A3").Value = 19
$ws.Range("B3").Value = "高雄市鼓山區龍北段00344000建號"
$ws.Range("C3").Value = 11097.76
$ws.Range("D3").Value = "100000分之1107"
$ws.Range("E3").Value = "吳美惠"
$ws.Range("F3").Value = "102年07月22曰"
$ws.Range("G3").Value = "買賣"
$ws.Range("H3").Value = "1020000(共同使用部分含2車位）"

$ws = $wb.Worksheets.Item(3)
$ws.Range("B1").Value = "廠牌型號"
$ws.Range("C1").Value = "汽缸容量"
$ws.Range("D1").Value = "所有人"
$ws.Range("E1").Value = "登記（取得）時間"
$ws.Range("F1").Value = "登記（取得）原因"
$ws.Range("G1").Value = "取得價額"
$ws.Range("A2").Value = 29
$ws.Range("B2").Value = "BluebirdSylphyGllTS(客車)"
$ws.Range("C2").Value = 1997
$ws.Range("D2").Value = "吳美惠"
$ws.Range("E2").Value = "98年10月12日"
$ws.Range("F2").Value = "買賣"
$ws.Range("G2").Value = 740000
$ws.Range("A3").Value = 30
$ws.Range("B3").Value = "納智捷L91MLD"
$ws.Range("C3").Value = 2198
$ws.Range("D3").Value = "林國正"
$ws.Range("E3").Value = "102年01月15曰"
$ws.Range("F3").Value = "買賣."
$ws.Range("G3").Value = 930000

$ws = $wb.Worksheets.Item(4)
$ws.Range("B1").Value = "幣別"
$ws.Range("C1").Value = "所有人"
$ws.Range("D1").Value = "新臺幣總額或折合新臺幣總額"
$ws.Range("A2").Value = 40
$ws.Range("B2").Value = "新臺幣"
$ws.Range("C2").Value = "林國正"
$ws.Range("D2").Value = 21000
$ws.Range("A3").Value = 41
$ws.Range("B3").Value = "新臺幣"
$ws.Range("C3").Value = "吳美惠"
$ws.Range("D3").Value = 16000

$ws = $wb.Worksheets.Item(5)
$ws.Range("B1").Value = "存放機構(應敘明分支機構）"
$ws.Range("C1").Value = "種類"
$ws.Range("D1").Value = "幣別"
$ws.Range("E1").Value = "所有人"
$ws.Range("F1").Value = "新臺幣總額或折合新臺幣總額"
$ws.Range("A2").Value = 46
$ws.Range("B2").Value = "國泰世華商業銀行南高雄分行"
$ws.Range("C2").Value = "活期存款"
$ws.Range("D2").Value = "新臺幣"
$ws.Range("E2").Value = "吳美惠"
$ws.Range("F2").Value = 3014
$ws.Range("A3").Value = 47
$ws.Range("B3").Value = "台北富邦商業銀行基隆路分行"
$ws.Range("C3").Value = "活期存款"
$ws.Range("D3").Value = "新臺幣"
$ws.Range("E3").Value = "林國正"
$ws.Range("F3").Value = 720
$ws.Range("A4").Value = 48
$ws.Range("B4").Value = "合作金庫商業銀行港都分行"
$ws.Range("C4").Value = "活期存款"
$ws.Range("D4").Value = "新臺幣"
$ws.Range("E4").Value = "林國正"
$ws.Range("F4").Value = 790443
$ws.Range("A5").Value = 49
$ws.Range("B5").Value = "台北公館郵局（第13支局）"
$ws.Range("C5").Value = "活期存款"
$ws.Range("D5").Value = "新臺幣"
$ws.Range("E5").Value = "林國正"
$ws.Range("F5").Value = 1042
$ws.Range("A6").Value = 50
$ws.Range("B6").Value = "梧棲郵局（第42支局）"
$ws.Range("C6").Value = "活期存款"
$ws.Range("D6").Value = "新臺幣"
$ws.Range("E6").Value = "林國正"
$ws.Range("F6").Value = 380
$ws.Range("A7").Value = 51
$ws.Range("B7").Value = "高雄金福路郵局(第44支局）"
$ws.Range("C7").Value = "中華郵政劃撥儲金"
$ws.Range("D7").Value = "新臺幣"
$ws.Range("E7").Value = "林國正"
$ws.Range("F7").Value = 27
$ws.Range("A8").Value = 52
$ws.Range("B8").Value = "華南商業銀行北高雄分行"
$ws.Range("C8").Value = "活期存款"
$ws.Range("D8").Value = "新臺幣"
$ws.Range("E8").Value = "吳美惠"
$ws.Range("F8").Value = 61589
$ws.Range("A9").Value = 53
$ws.Range("B9").Value = "彰化商業銀行南高雄分行"
$ws.Range("C9").Value = "活期存款"
$ws.Range("D9").Value = "新臺幣"
$ws.Range("E9").Value = "林國正"
$ws.Range("F9").Value = 18785
$ws.Range("A10").Value = 54
$ws.Range("B10").Value = "高雄銀行七賢分行"
$ws.Range("C10").Value = "活期存款"
$ws.Range("D10").Value = "新臺幣"
$ws.Range("E10").Value = "林國正"
$ws.Range("F10").Value = 1521915
$ws.Range("A11").Value = 55
$ws.Range("B11").Value = "高雄銀行營業部"
$ws.Range("C11").Value = "活期存款"
$ws.Range("D11").Value = "新臺幣"
$ws.Range("E11").Value = "吳美惠"
$ws.Range("F11").Value = 84971
$ws.Range("A12").Value = 56
$ws.Range("B12").Value = "兆豐國際商業銀行三民分行"
$ws.Range("C12").Value = "活期儲蓄存款"
$ws.Range("D12").Value = "新臺幣"
$ws.Range("E12").Value = "吳美惠"
$ws.Range("F12").Value = 372270

$ws = $wb.Worksheets.Item(6)
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "owner"
$ws.Range("D1").Value = "quantity"
$ws.Range("E1").Value = "face_value"
$ws.Range("F1").Value = "currency"
$ws.Range("G1").Value = "total"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"
$ws.Range("A2").Value = 63
$ws.Range("B2").Value = "潤泰全球股份有限公司"
$ws.Range("C2").Value = "吳美惠"
$ws.Range("D2").Value = 423
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = "新臺幣"
$ws.Range("G2").Value = 4230
$ws.Range("H2").Value = "stock"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2013-12-24"
$ws.Range("K2").Value = "林國正"
$ws.Range("L2").Value = 1742
$ws.Range("M2").Value = "tmp399c1"
$ws.Range("N2").Value = 63

$ws = $wb.Worksheets.Item(7)
$ws.Range("B1").Value = "種類"
$ws.Range("C1").Value = "債務人"
$ws.Range("D1").Value = "債權人及地址"
$ws.Range("E1").Value = "餘額"
$ws.Range("F1").Value = "取得（發生）時間"
$ws.Range("G1").Value = "取得（發生）原因"
$ws.Range("A2").Value = 101
$ws.Range("B2").Value = "房屋貸款"
$ws.Range("C2").Value = "吳美惠"
$ws.Range("D2").Value = "兆豐國際商業銀行三民分行高雄市鼓山區中華一路225號"
$ws.Range("E2").Value = 5890212
$ws.Range("F2").Value = "102年07月23日"
$ws.Range("G2").Value = "買房屋貸款"
